$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf = New-Object 'object[,]' 24,5
$bf[0,0] = 0.3069574617624653
$bf[0,1] = 0.04707455330056121
$bf[0,2] = 0.03034126822065275
$bf[0,3] = 0.1597681957234514
$bf[0,4] = 0.7466487791299841
$bf[1,0] = 0.2734750974589133
$bf[1,1] = 0.04344008859901294
$bf[1,2] = 0.02844548043141515
$bf[1,3] = 0.1487329861242443
$bf[1,4] = 0.7437578374085021
$bf[2,0] = 0.252954556317178
$bf[2,1] = 0.04119035834491314
$bf[2,2] = 0.02726995512389152
$bf[2,3] = 0.1420638794958577
$bf[2,4] = 0.7424312632612455
$bf[3,0] = 0.2446021479260594
$bf[3,1] = 0.04026904427728084
$bf[3,2] = 0.02678805104179105
$bf[3,3] = 0.1393727600291967
$bf[3,4] = 0.7420034464434622
$bf[4,0] = 0.2432158447921324
$bf[4,1] = 0.04011578814749583
$bf[4,2] = 0.02670785874732928
$bf[4,3] = 0.1389275021947896
$bf[4,4] = 0.7419392190108951
$bf[5,0] = 0.2528418723109667
$bf[5,1] = 0.04117795147899272
$bf[5,2] = 0.0272634675764607
$bf[5,3] = 0.1420274786953968
$bf[5,4] = 0.7424250369616701
$bf[6,0] = 0.2954051231704398
$bf[6,1] = 0.04582517733172153
$bf[6,2] = 0.029690002582754
$bf[6,3] = 0.1559410164013855
$bf[6,4] = 0.7455588981920727
$bf[7,0] = 0.3791585438694369
$bf[7,1] = 0.05479322258470631
$bf[7,2] = 0.03435626375856771
$bf[7,3] = 0.184081356257586
$bf[7,4] = 0.7552642961534772
$bf[8,0] = 0.4408562900354127
$bf[8,1] = 0.06129266680953549
$bf[8,2] = 0.03772748102181822
$bf[8,3] = 0.2052944720597552
$bf[8,4] = 0.764569173393312
$bf[9,0] = 0.4689578846499671
$bf[9,1] = 0.06422987158329363
$bf[9,2] = 0.03924856618433381
$bf[9,3] = 0.215065472280827
$bf[9,4] = 0.7692753502321921
$bf[10,0] = 0.4796039507425007
$bf[10,1] = 0.0653392944196014
$bf[10,2] = 0.03982274263158558
$bf[10,3] = 0.21878315229376
$bf[10,4] = 0.7711255631687237
$bf[11,0] = 0.4773109320678373
$bf[11,1] = 0.06510048698271476
$bf[11,2] = 0.03969916509885962
$bf[11,3] = 0.2179816966768229
$bf[11,4] = 0.7707240582584944
$bf[12,0] = 0.4698336531352822
$bf[12,1] = 0.06432120146916986
$bf[12,2] = 0.03929584077430093
$bf[12,3] = 0.2153709738461913
$bf[12,4] = 0.7694262037092443
$bf[13,0] = 0.4652541905035719
$bf[13,1] = 0.06384349649299281
$bf[13,2] = 0.03904855440776345
$bf[13,3] = 0.2137741307960255
$bf[13,4] = 0.7686400976983236
$bf[14,0] = 0.4390204520542795
$bf[14,1] = 0.06110031930472815
$bf[14,2] = 0.03762782061422598
$bf[14,3] = 0.2046583677667755
$bf[14,4] = 0.7642711402333759
$bf[15,0] = 0.4229355581424556
$bf[15,1] = 0.05941246556679403
$bf[15,2] = 0.03675302520878176
$bf[15,3] = 0.1990972934537965
$bf[15,4] = 0.7617121748877835
$bf[16,0] = 0.4136872707346981
$bf[16,1] = 0.05843983072834646
$bf[16,2] = 0.03624869130064923
$bf[16,3] = 0.1959100917404868
$bf[16,4] = 0.7602848776923921
$bf[17,0] = 0.4105565424548274
$bf[17,1] = 0.05811020087192276
$bf[17,2] = 0.03607773160440786
$bf[17,3] = 0.1948329083922076
$bf[17,4] = 0.7598092706155768
$bf[18,0] = 0.4246474822432447
$bf[18,1] = 0.05959232975158102
$bf[18,2] = 0.03684627045753786
$bf[18,3] = 0.1996880998977488
$bf[18,4] = 0.7619799700057328
$bf[19,0] = 0.4720297897237913
$bf[19,1] = 0.06455017371138183
$bf[19,2] = 0.03941435667330495
$bf[19,3] = 0.2161373267239881
$bf[19,4] = 0.7698055670632016
$bf[20,0] = 0.503023492289941
$bf[20,1] = 0.06777388091623493
$bf[20,2] = 0.04108209599903745
$bf[20,3] = 0.2269906132914343
$bf[20,4] = 0.7753169126578996
$bf[21,0] = 0.4864792752467224
$bf[21,1] = 0.0660548533064258
$bf[21,2] = 0.04019297625615792
$bf[21,3] = 0.2211885382546726
$bf[21,4] = 0.7723390841501896
$bf[22,0] = 0.4238735240286928
$bf[22,1] = 0.05951102020779331
$bf[22,2] = 0.03680411864618094
$bf[22,3] = 0.1994209653795096
$bf[22,4] = 0.7618587631313645
$bf[23,0] = 0.3564714194400551
$bf[23,1] = 0.05238276043010615
$bf[23,2] = 0.03310387632460277
$bf[23,3] = 0.1763752981699653
$bf[23,4] = 0.752257190075035
$ws.Range("B2:F25").Value = $bf

$k = New-Object 'object[,]' 24,1
$k[0,0] = 0.2775569654750143
$k[1,0] = 0.243139811749856
$k[2,0] = 0.2219865104311083
$k[3,0] = 0.2133614458601585
$k[4,0] = 0.2119289739284369
$k[5,0] = 0.2218702092406346
$k[6,0] = 0.2656945045555119
$k[7,0] = 0.3514550262920011
$k[8,0] = 0.4143443601893466
$k[9,0] = 0.4429267800225034
$k[10,0] = 0.4537461332288331
$k[11,0] = 0.4514161848230458
$k[12,0] = 0.4438169810439661
$k[13,0] = 0.4391616915108614
$k[14,0] = 0.41247586819307
$k[15,0] = 0.3960979615617362
$k[16,0] = 0.3866753764795305
$k[17,0] = 0.3834846435057671
$k[18,0] = 0.397841673867191
$k[19,0] = 0.4460491683835244
$k[20,0] = 0.4775309447403515
$k[21,0] = 0.4607309166050868
$k[22,0] = 0.3970533623996459
$k[23,0] = 0.328274698999337
$ws.Range("K2:K25").Value = $k

$mo = New-Object 'object[,]' 24,3
$mo[0,0] = 0.2255241579484561
$mo[0,1] = 1.501335007875126
$mo[0,2] = 2.576247709650858
$mo[1,0] = 0.2033048863317291
$mo[1,1] = 1.518021207761241
$mo[1,2] = 2.579888993381729
$mo[2,0] = 0.1897423394864575
$mo[2,1] = 1.528791050497801
$mo[2,2] = 2.583585944296644
$mo[3,0] = 0.1842356117007355
$mo[3,1] = 1.533311674359807
$mo[3,2] = 2.585459812436724
$mo[4,0] = 0.1833224377961216
$mo[4,1] = 1.534070280978531
$mo[4,2] = 2.585793149894755
$mo[5,0] = 0.1896679924169007
$mo[5,1] = 1.528851483525036
$mo[5,2] = 2.58360972876261
$mo[6,0] = 0.2178462812302371
$mo[6,1] = 1.506979530991161
$mo[6,2] = 2.577199889859344
$mo[7,0] = 0.2737447342510677
$mo[7,1] = 1.468251585543371
$mo[7,2] = 2.576231566654911
$mo[8,0] = 0.3152148378202497
$mo[8,1] = 1.442337222943445
$mo[8,2] = 2.582607475424311
$mo[9,0] = 0.3341704033445154
$mo[9,1] = 1.431099715214041
$mo[9,2] = 2.587050496684213
$mo[10,0] = 0.3413615212004188
$mo[10,1] = 1.426923668676894
$mo[10,2] = 2.588954982172083
$mo[11,0] = 0.3398122041824365
$mo[11,1] = 1.42781952394005
$mo[11,2] = 2.588534940221422
$mo[12,0] = 0.33476175911116
$mo[12,1] = 1.430754558338892
$mo[12,2] = 2.587202729373217
$mo[13,0] = 0.3316699168271384
$mo[13,1] = 1.432562689200186
$mo[13,2] = 2.586415629542131
$mo[14,0] = 0.3139778726423188
$mo[14,1] = 1.443082718797049
$mo[14,2] = 2.582348170388258
$mo[15,0] = 0.3031476030184095
$mo[15,1] = 1.449677635650861
$mo[15,2] = 2.580248143317363
$mo[16,0] = 0.2969268451991098
$mo[16,1] = 1.453522719494321
$mo[16,2] = 2.579185447808214
$mo[17,0] = 0.2948220671307382
$mo[17,1] = 1.454833502841048
$mo[17,2] = 2.578850567240721
$mo[18,0] = 0.3042996210159075
$mo[18,1] = 1.44897022716065
$mo[18,2] = 2.580456667367258
$mo[19,0] = 0.3362448426873215
$mo[19,1] = 1.429890313052695
$mo[19,2] = 2.587588005219203
$mo[20,0] = 0.3571990076318841
$mo[20,1] = 1.417883074047749
$mo[20,2] = 2.593542876277354
$mo[21,0] = 0.3460084019974943
$mo[21,1] = 1.424249199263469
$mo[21,2] = 2.590246174757908
$mo[22,0] = 0.3037787759033534
$mo[22,1] = 1.449289879649237
$mo[22,2] = 2.580361943113985
$mo[23,0] = 0.2585529892774545
$mo[23,1] = 1.478282736382013
$mo[23,2] = 2.575249986083747
$ws.Range("M2:O25").Value = $mo
